# Apply the "Update countries & provincias Spain" edit:
#  - Update the "Datos actualizados" timestamp string (A1)
#  - Update the Murcia row (row 31) figures: Casos totales, Casos activos, Recuperados, Muertes

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 21:22"

# Update Murcia's row (row 31): B31, C31, D31, E31
$ws.Range("B31").Value = 1463
$ws.Range("C31").Value = 353
$ws.Range("D31").Value = 1009
$ws.Range("E31").Value = 101
